$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Find.Execute("2025-07-11 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-12 Saturday", 2) | Out-Null

# Update every answer cell in the table by its fixed (row, column) position
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "91-61=30"
$t.Cell(1,2).Range.Text = "61-9=52"
$t.Cell(1,3).Range.Text = "36+37=73"
$t.Cell(1,4).Range.Text = "0+58=58"
$t.Cell(1,5).Range.Text = "77-22=55"
$t.Cell(2,1).Range.Text = "31+59=90"
$t.Cell(2,2).Range.Text = "77-6=71"
$t.Cell(2,3).Range.Text = "2+72=74"
$t.Cell(2,4).Range.Text = "4+28=32"
$t.Cell(2,5).Range.Text = "77-65=12"
$t.Cell(3,1).Range.Text = "19+76=95"
$t.Cell(3,2).Range.Text = "22+19=41"
$t.Cell(3,3).Range.Text = "92-34=58"
$t.Cell(3,4).Range.Text = "72-17=55"
$t.Cell(3,5).Range.Text = "49+41=90"
$t.Cell(4,1).Range.Text = "81-74=7"
$t.Cell(4,2).Range.Text = "77-57=20"
$t.Cell(4,3).Range.Text = "34+54=88"
$t.Cell(4,4).Range.Text = "99-93=6"
$t.Cell(4,5).Range.Text = "93-53=40"
$t.Cell(5,1).Range.Text = "48-15=33"
$t.Cell(5,2).Range.Text = "29+58=87"
$t.Cell(5,3).Range.Text = "30+23=53"
$t.Cell(5,4).Range.Text = "53+4=57"
$t.Cell(5,5).Range.Text = "97-33=64"
$t.Cell(6,1).Range.Text = "86-57=29"
$t.Cell(6,2).Range.Text = "96-1=95"
$t.Cell(6,3).Range.Text = "43-42=1"
$t.Cell(6,4).Range.Text = "88-47=41"
$t.Cell(6,5).Range.Text = "93-72=21"
$t.Cell(7,1).Range.Text = "40-5=35"
$t.Cell(7,2).Range.Text = "69-4=65"
$t.Cell(7,3).Range.Text = "99-94=5"
$t.Cell(7,4).Range.Text = "27+26=53"
$t.Cell(7,5).Range.Text = "5+10=15"
$t.Cell(8,1).Range.Text = "31+12=43"
$t.Cell(8,2).Range.Text = "59-6=53"
$t.Cell(8,3).Range.Text = "83-76=7"
$t.Cell(8,4).Range.Text = "71+16=87"
$t.Cell(8,5).Range.Text = "31-20=11"
$t.Cell(9,1).Range.Text = "82-16=66"
$t.Cell(9,2).Range.Text = "78-58=20"
$t.Cell(9,3).Range.Text = "91-80=11"
$t.Cell(9,4).Range.Text = "64-10=54"
$t.Cell(9,5).Range.Text = "51+44=95"
$t.Cell(10,1).Range.Text = "86-4=82"
$t.Cell(10,2).Range.Text = "43-30=13"
$t.Cell(10,3).Range.Text = "89-52=37"
$t.Cell(10,4).Range.Text = "22-19=3"
$t.Cell(10,5).Range.Text = "15+37=52"
$t.Cell(11,1).Range.Text = "99-81=18"
$t.Cell(11,2).Range.Text = "56-15=41"
$t.Cell(11,3).Range.Text = "20+54=74"
$t.Cell(11,4).Range.Text = "90-51=39"
$t.Cell(11,5).Range.Text = "74-9=65"
$t.Cell(12,1).Range.Text = "93-26=67"
$t.Cell(12,2).Range.Text = "5-3=2"
$t.Cell(12,3).Range.Text = "38+44=82"
$t.Cell(12,4).Range.Text = "24-22=2"
$t.Cell(12,5).Range.Text = "62-22=40"
$t.Cell(13,1).Range.Text = "64-11=53"
$t.Cell(13,2).Range.Text = "33-24=9"
$t.Cell(13,3).Range.Text = "9+87=96"
$t.Cell(13,4).Range.Text = "97-68=29"
$t.Cell(13,5).Range.Text = "37+22=59"
$t.Cell(14,1).Range.Text = "71-20=51"
$t.Cell(14,2).Range.Text = "68+2=70"
$t.Cell(14,3).Range.Text = "96-61=35"
$t.Cell(14,4).Range.Text = "32+20=52"
$t.Cell(14,5).Range.Text = "23+65=88"
$t.Cell(15,1).Range.Text = "14+85=99"
$t.Cell(15,2).Range.Text = "3+17=20"
$t.Cell(15,3).Range.Text = "16+38=54"
$t.Cell(15,4).Range.Text = "62-49=13"
$t.Cell(15,5).Range.Text = "27-13=14"
$t.Cell(16,1).Range.Text = "38-13=25"
$t.Cell(16,2).Range.Text = "39+50=89"
$t.Cell(16,3).Range.Text = "58-7=51"
$t.Cell(16,4).Range.Text = "97-62=35"
$t.Cell(16,5).Range.Text = "36+39=75"
$t.Cell(17,1).Range.Text = "54-8=46"
$t.Cell(17,2).Range.Text = "99-27=72"
$t.Cell(17,3).Range.Text = "96-56=40"
$t.Cell(17,4).Range.Text = "92-66=26"
$t.Cell(17,5).Range.Text = "38-9=29"
$t.Cell(18,1).Range.Text = "29+57=86"
$t.Cell(18,2).Range.Text = "16-5=11"
$t.Cell(18,3).Range.Text = "68+19=87"
$t.Cell(18,4).Range.Text = "0+62=62"
$t.Cell(18,5).Range.Text = "13+9=22"
$t.Cell(19,1).Range.Text = "43+13=56"
$t.Cell(19,2).Range.Text = "40+58=98"
$t.Cell(19,3).Range.Text = "24+17=41"
$t.Cell(19,4).Range.Text = "45-0=45"
$t.Cell(19,5).Range.Text = "67-21=46"
$t.Cell(20,1).Range.Text = "85-58=27"
$t.Cell(20,2).Range.Text = "49+1=50"
$t.Cell(20,3).Range.Text = "95-83=12"
$t.Cell(20,4).Range.Text = "2+15=17"
$t.Cell(20,5).Range.Text = "49+40=89"
